# Automatische test-sync: 2025-08-05 18:46:50
$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 40 with the latest test mail entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A40").Value = "Bel jij klant Jansen even?"
$logs.Range("B40").Value = "mailmind.test@zohomail.eu"
$logs.Range("C40").Value = "Testmail #19: Bel jij klant Jansen even?"
$logs.Range("D40").Value = "Klantenservice / Contact"
$logs.Range("E40").Value = "Bedankt, we hebben dit doorgestuurd naar klantenservice@bedrijf.nl."
$logs.Range("F40").Value = "2025-08-05 18:45:56"
$logs.Range("G40").Value = "Ja"
$logs.Range("H40").Value = "Ja"
$logs.Range("I40").Value = "Nee"
$logs.Range("J40").Value = "Nee"

# --- Expand the conditional formatting ranges to include the new row 40 ---
$logs.Range("D2:D39").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D40"))
$logs.Range("G2:G39").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G40"))
$logs.Range("H2:H39").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H40"))
$logs.Range("I2:I39").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I40"))
$logs.Range("J2:J39").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J40"))

# --- Sheet "Dashboard": category counts changed order/values ---
# "Klantenservice / Contact" now ties "Inkoop / Bestellingen" at 5 and moves up to row 3
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Klantenservice / Contact"
$dash.Range("B3").Value = 5
$dash.Range("A4").Value = "Inkoop / Bestellingen"
$dash.Range("B4").Value = 5
